$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("SigninData")

# H2: formula now references SigninData!I2 (which will contain "Testing123")
$ws.Range("H2").Formula = "=SigninData!I2"

# I2: becomes the literal text "Testing123" (stored as shared string)
$ws.Range("I2").Value = "Testing123"

# J2: attempt count incremented from 10 to 11
$ws.Range("J2").Value = 11
